# Update cryptocurrency price/volume data per the Sun Jun 18 22:41:11 UTC 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its exact textual representation (e.g. trailing
# zeros like "1.000") instead of being reinterpreted as a number by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "26.437.25"
$ws.Range("E2").Value = "  -0.37%  "

# Row 3
$ws.Range("D3").Value = "1.726.65"
$ws.Range("E3").Value = "  -0.19%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").Value = "243.00"
$ws.Range("E5").Value = "  -1.02%  "

# Row 6
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.06%  "

# Row 7
$ws.Range("D7").Value = "0.4859"
$ws.Range("E7").Value = "  +1.06%  "

# Row 8
$ws.Range("D8").Value = "0.2613"
$ws.Range("E8").Value = "  -2.27%  "

# Row 9
$ws.Range("D9").Value = "0.06192"
$ws.Range("E9").Value = "  -0.52%  "

# Row 10
$ws.Range("D10").Value = "1.730.66"
$ws.Range("E10").Value = "  +0.04%  "

# Row 11
$ws.Range("D11").Value = "0.07013"
$ws.Range("E11").Value = "  -2.11%  "

# Row 12
$ws.Range("D12").Value = "15.44"
$ws.Range("E12").Value = "  -1.73%  "

# Row 13
$ws.Range("D13").Value = "4.546"
$ws.Range("E13").Value = "  +0.10%  "

# Row 14
$ws.Range("D14").Value = "0.5989"
$ws.Range("E14").Value = "  -2.90%  "

# Row 15
$ws.Range("D15").Value = "77.29"
$ws.Range("E15").Value = "  +0.06%  "

# Row 16
$ws.Range("E16").Value = "  +0.07%  "

# Row 17
$ws.Range("D17").Value = "26.452.50"
$ws.Range("E17").Value = "  -0.33%  "

# Row 18
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.04%  "

# Row 19
$ws.Range("D19").Value = "0.000007199"
$ws.Range("E19").Value = "  +3.52%  "

# Row 20
$ws.Range("D20").Value = "11.41"
$ws.Range("E20").Value = "  -2.15%  "

# Row 21
$ws.Range("D21").Value = "1.954.47"
$ws.Range("E21").Value = "  +0.24%  "

# Row 22
$ws.Range("D22").Value = "4.493"
$ws.Range("E22").Value = "  -0.90%  "

# Row 23
$ws.Range("D23").Value = "8.585"
$ws.Range("E23").Value = "  -3.86%  "

# Row 24
$ws.Range("D24").Value = "5.179"
$ws.Range("E24").Value = "  -2.09%  "

# Row 25
$ws.Range("D25").Value = "138.27"
$ws.Range("E25").Value = "  +1.23%  "

# Row 26
$ws.Range("D26").Value = "15.25"
$ws.Range("E26").Value = "  -0.70%  "

# Row 27
$ws.Range("D27").Value = "1.410"
$ws.Range("E27").Value = "  +0.32%  "

# Row 28
$ws.Range("D28").Value = "106.97"
$ws.Range("E28").Value = "  -0.16%  "

# Row 29
$ws.Range("D29").Value = "1.720"
$ws.Range("E29").Value = "  -4.32%  "

# Row 30
$ws.Range("D30").Value = "3.956"
$ws.Range("E30").Value = "  -0.97%  "

# Row 31
$ws.Range("D31").Value = "0.07951"
$ws.Range("E31").Value = "  -0.62%  "

# Row 32
$ws.Range("D32").Value = "3.689"
$ws.Range("E32").Value = "  -0.91%  "

# Row 33
$ws.Range("D33").Value = "0.04517"
$ws.Range("E33").Value = "  -1.62%  "

# Row 34
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "2.614"
$ws.Range("E34").Value = "  -0.13%  "

# Row 35
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "0.9994"
$ws.Range("E35").Value = "  +0.30%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.6236"
$ws.Range("E36").Value = "  -2.16%  "

# Row 37
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "0.9075"
$ws.Range("E37").Value = "  -1.78%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "1.992"
$ws.Range("E38").Value = "  -4.89%  "

# Row 39
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.406"
$ws.Range("E39").Value = "  -0.11%  "

# Row 40
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").Value = "1.001"
$ws.Range("E40").Value = "  -0.56%  "

# Row 41
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.01489"
$ws.Range("E41").Value = "  -1.20%  "

# Row 42
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "100.23"
$ws.Range("E42").Value = "  -4.40%  "

# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "5.398"
$ws.Range("E43").Value = "  -3.51%  "

# Row 44
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.3861"
$ws.Range("E44").Value = "  -1.09%  "

# Row 45
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "6.680"
$ws.Range("E45").Value = "  -4.30%  "

# Row 46
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "0.1155"
$ws.Range("E46").Value = "  -2.51%  "

# Row 47
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.05360"
$ws.Range("E47").Value = "  +0.38%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "7.736"
$ws.Range("E48").Value = "  -1.59%  "

# Row 49
$ws.Range("D49").Value = "30.22"
$ws.Range("E49").Value = "  -2.58%  "

# Row 50
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "1.253"
$ws.Range("E50").Value = "  -1.11%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "51.03"
$ws.Range("E51").Value = "  -0.56%  "
